# Add more expenses: append rows 1825-1848 to the "total" sheet's expense log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: the new transaction row plus the four "derived" columns (A-D, J, K).
# E:I (month/year/weekday-number/weekday-name/month-abbr) are filled afterwards
# via formulas, exactly like the rest of the sheet.
$rows = @(
  @{ Row=1825; A=45414; B='food'; C='menu 2'; D=2.9;  DF=$null;            J='infineon';           K=$null },
  @{ Row=1826; A=45414; B='food'; C='dessert'; D=0.79; DF=$null;            J='infineon';           K=$null },
  @{ Row=1827; A=45415; B='food'; C='menu 1'; D=3.3;  DF=$null;            J='infineon';           K=$null },
  @{ Row=1828; A=45415; B='food'; C='dessert-creme'; D=0.79; DF=$null;     J='infineon';           K=$null },
  @{ Row=1829; A=45416; B='food'; C='flying goose sriracha wasabi 200ml'; D=2.6; DF=$null; J='modern asia market'; K='graz' },
  @{ Row=1830; A=45416; B='food'; C='pocky cookies & cream 40g'; D=2.5; DF=$null; J='modern asia market'; K='graz' },
  @{ Row=1831; A=45416; B='food'; C='basic kuchenrolle'; D=$null;         DF='=2.39/2';           J='billa';              K='graz' },
  @{ Row=1832; A=45416; B='food'; C='vanillejog. Schokomuesli'; D=$null;  DF='=3.82-1.91';        J='billa';              K='graz' },
  @{ Row=1833; A=45418; B='food'; C='salat gross'; D=2.4; DF=$null;       J='infineon';           K=$null },
  @{ Row=1834; A=45419; B='food'; C='grill'; D=4.16; DF=$null;            J='infineon';           K=$null },
  @{ Row=1835; A=45420; B='food'; C='menu 2'; D=2.9;  DF=$null;           J='infineon';           K=$null },
  @{ Row=1836; A=45419; B='food'; C='protein pita roll chicken'; D=3.79; DF=$null; J='billa';      K='villach' },
  @{ Row=1837; A=45419; B='food'; C='nuss schnecke'; D=1.69; DF=$null;   J='billa';              K='villach' },
  @{ Row=1838; A=45419; B='food'; C='franziskaner fi.'; D=1.29; DF=$null; J='billa';              K='villach' },
  @{ Row=1839; A=45420; B='food'; C='rm apfel rot tasse 0.797kgx1.99'; D=1.59; DF=$null; J='billa'; K='villach' },
  @{ Row=1840; A=45420; B='food'; C='nuss schnecke'; D=1.69; DF=$null;   J='billa';              K='villach' },
  @{ Row=1841; A=45416; B='gas';  C='gas'; D=$null;                      DF='=49.33/2';          J='avanti';             K='klagenfurt' },
  @{ Row=1842; A=45418; B='food'; C='apfel 0,998x1,89kg'; D=1.89; DF=$null; J='hofer';            K='klagenfurt' },
  @{ Row=1843; A=45418; B='food'; C='blattsalat-mix'; D=0.99; DF=$null;  J='hofer';              K='klagenfurt' },
  @{ Row=1844; A=45418; B='food'; C='rohkost-salatmix'; D=0.99; DF=$null; J='hofer';             K='klagenfurt' },
  @{ Row=1845; A=45420; B='food'; C='spar curly-fries 600g'; D=$null;    DF='=(2.99-1)/2';       J='eurospar';           K='klagenfurt' },
  @{ Row=1846; A=45420; B='food'; C='spar gitterpommes'; D=$null;        DF='=(2.99-1)/2';       J='eurospar';           K='klagenfurt' },
  @{ Row=1847; A=45420; B='food'; C='sbudget mischsalat'; D=1.99; DF=$null; J='eurospar';         K='klagenfurt' },
  @{ Row=1848; A=45420; B='food'; C='spar highprotknm 500g'; D=3.19; DF=$null; J='eurospar';      K='klagenfurt' }
)

foreach ($r in $rows) {
  $ws.Cells.Item($r.Row, 1).Value = $r.A
  $ws.Cells.Item($r.Row, 2).Value = $r.B
  $ws.Cells.Item($r.Row, 3).Value = $r.C
  if ($r.DF -ne $null) {
    $ws.Cells.Item($r.Row, 4).Formula = $r.DF
  } else {
    $ws.Cells.Item($r.Row, 4).Value = $r.D
  }
  $ws.Cells.Item($r.Row, 10).Value = $r.J
  if ($r.K -ne $null) {
    $ws.Cells.Item($r.Row, 11).Value = $r.K
  }
}

# Columns E-I are the same formulas used throughout the sheet, filled down
# across the whole new block (MONTH / YEAR / WEEKDAY / weekday-name / month-abbr).
$ws.Range("E1825:E1848").Formula = "=MONTH(A1825)"
$ws.Range("F1825:F1848").Formula = "=YEAR(A1825)"
$ws.Range("G1825:G1848").Formula = "=WEEKDAY(A1825, 2)"
$ws.Range("H1825:H1848").Formula = "=CHOOSE(WEEKDAY(A1825, 2), ""Monday"", ""Tuesday"",""Wednesday"", ""Thursday"", ""Friday"", ""Saturday"",""Sunday"")"
$ws.Range("I1825:I1848").Formula = "=TEXT(A1825, ""MMM"")"

# Keep the frozen-pane view / selection in sync with the newly added rows,
# matching how Excel leaves the sheet after scrolling to the bottom and
# selecting the last-edited cell.
$ws.Application.ActiveWindow.ScrollRow = 1826
$ws.Range("H1835").Select()
